$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-80 down to 9-81
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data record
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 45063
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112035
$ws.Range("G8").Value = "Bruselas (repollito)"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 52
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21000
$ws.Range("N8").Value = "`$/malla 15 kilos"
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 1400
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = "Hortaliza"
